$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 166667280
$ws.Range("I107").Value = 333333600
$ws.Range("J107").Value = 965.3333
$ws.Range("K107").Value = 333333600
$ws.Range("L107").Value = 965.3333
$ws.Range("M107").Value = -333331680
$ws.Range("N107").Value = -4805.3333
$ws.Range("H111").Value = 2027.6666
$ws.Range("I111").Value = 1730.6154
$ws.Range("J111").Value = 2800
$ws.Range("K111").Value = 5191.8462
$ws.Range("L111").Value = 8400
$ws.Range("M111").Value = -2124.8462
$ws.Range("N111").Value = -14534
$ws.Range("H112").Value = 1714.5
$ws.Range("I112").Value = 2149.5
$ws.Range("J112").Value = 1590.2142
$ws.Range("K112").Value = 6448.5
$ws.Range("L112").Value = 4770.642599999999
$ws.Range("M112").Value = -5340.5
$ws.Range("N112").Value = -6986.642599999999
$ws.Range("H115").Value = 349.5
$ws.Range("I115").Value = 349.5
$ws.Range("K115").Value = 1048.5
$ws.Range("M115").Value = 518.5
$ws.Range("H116").Value = 2166.795
$ws.Range("I116").Value = 2163
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 2163
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = 1279
$ws.Range("N116").Value = -9084
$ws.Range("H118").Value = 4207.8076
$ws.Range("I118").Value = 396.8
$ws.Range("J118").Value = 9404.637000000001
$ws.Range("K118").Value = 1190.4
$ws.Range("L118").Value = 28213.911
$ws.Range("M118").Value = 466.5999999999999
$ws.Range("N118").Value = -31527.911
$ws.Range("H129").Value = 1768.9302
$ws.Range("J129").Value = 2099.7144
$ws.Range("L129").Value = 6299.1432
$ws.Range("N129").Value = -16299.1432
$ws.Range("H137").Value = 1604.421
$ws.Range("I137").Value = 1552.4166
$ws.Range("J137").Value = 1693.5714
$ws.Range("K137").Value = 4657.2498
$ws.Range("L137").Value = 5080.7142
$ws.Range("M137").Value = -2107.2498
$ws.Range("N137").Value = -10180.7142
$ws.Range("H140").Value = 36850
$ws.Range("J140").Value = 36850
$ws.Range("L140").Value = 36850
$ws.Range("N140").Value = -47210

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1008.5833
$ws.Range("I2").Value = 941.8889
$ws.Range("J2").Value = 1208.6666
$ws.Range("K2").Value = 941.8889
$ws.Range("L2").Value = 1208.6666
$ws.Range("M2").Value = -828.8889
$ws.Range("N2").Value = -1434.6666
$ws.Range("H88").Value = 76925670
$ws.Range("I88").Value = 1566.6666
$ws.Range("J88").Value = 100002904
$ws.Range("K88").Value = 1566.6666
$ws.Range("L88").Value = 100002904
$ws.Range("M88").Value = -1160.6666
$ws.Range("N88").Value = -100003716
$ws.Range("H91").Value = 76925670
$ws.Range("I91").Value = 1566.6666
$ws.Range("J91").Value = 100002904
$ws.Range("K91").Value = 1566.6666
$ws.Range("L91").Value = 100002904
$ws.Range("M91").Value = -162.6666
$ws.Range("N91").Value = -100005712
$ws.Range("H110").Value = 7970.0625
$ws.Range("I110").Value = 8608.643
$ws.Range("J110").Value = 3500
$ws.Range("K110").Value = 8608.643
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = -6563.643
$ws.Range("N110").Value = -7590
$ws.Range("H116").Value = 1008.5833
$ws.Range("I116").Value = 941.8889
$ws.Range("J116").Value = 1208.6666
$ws.Range("K116").Value = 941.8889
$ws.Range("L116").Value = 1208.6666
$ws.Range("M116").Value = 1352.1111
$ws.Range("N116").Value = -5796.6666

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1008.5833
$ws.Range("I3").Value = 941.8889
$ws.Range("J3").Value = 1208.6666
$ws.Range("K3").Value = 941.8889
$ws.Range("L3").Value = 1208.6666
$ws.Range("M3").Value = -827.8889
$ws.Range("N3").Value = -1436.6666
$ws.Range("H86").Value = 8335305
$ws.Range("I86").Value = 16668810
$ws.Range("J86").Value = 1800
$ws.Range("K86").Value = 16668810
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -16667687
$ws.Range("N86").Value = -4046
$ws.Range("H89").Value = 8335305
$ws.Range("I89").Value = 16668810
$ws.Range("J89").Value = 1800
$ws.Range("K89").Value = 83344050
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -83338434
$ws.Range("N89").Value = -20232
$ws.Range("H105").Value = 2709
$ws.Range("I105").Value = 2155.7144
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 2155.7144
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -408.7143999999998
$ws.Range("N105").Value = -7494
$ws.Range("H107").Value = 2571.25
$ws.Range("I107").Value = 2370
$ws.Range("K107").Value = 2370
$ws.Range("M107").Value = -450

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10028734
$ws.Range("I31").Value = 8041115.5
$ws.Range("J31").Value = 13341431
$ws.Range("K31").Value = 8041115.5
$ws.Range("L31").Value = 13341431
$ws.Range("M31").Value = -8040820.5
$ws.Range("N31").Value = -13342021
$ws.Range("H34").Value = 10028734
$ws.Range("I34").Value = 8041115.5
$ws.Range("J34").Value = 13341431
$ws.Range("K34").Value = 8041115.5
$ws.Range("L34").Value = 13341431
$ws.Range("M34").Value = -8040913.5
$ws.Range("N34").Value = -13341835
$ws.Range("H107").Value = 20834424
$ws.Range("I107").Value = 33334230
$ws.Range("J107").Value = 1415.5555
$ws.Range("K107").Value = 33334230
$ws.Range("L107").Value = 1415.5555
$ws.Range("M107").Value = -33332310
$ws.Range("N107").Value = -5255.5555

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1216.9412
$ws.Range("I5").Value = 1199.2
$ws.Range("K5").Value = 3597.6
$ws.Range("M5").Value = -3485.6
$ws.Range("H122").Value = 679.8570999999999
$ws.Range("I122").Value = 393.22223
$ws.Range("J122").Value = 1195.8
$ws.Range("K122").Value = 3539.00007
$ws.Range("L122").Value = 10762.2
$ws.Range("M122").Value = -1089.00007
$ws.Range("N122").Value = -15662.2
$ws.Range("H135").Value = 1216.9412
$ws.Range("I135").Value = 1199.2
$ws.Range("K135").Value = 10792.8
$ws.Range("M135").Value = -8257.800000000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2442.2
$ws.Range("I113").Value = 1011
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 1011
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = 1159
$ws.Range("N113").Value = -7140
$ws.Range("H141").Value = 35164.5
$ws.Range("J141").Value = 35164.5
$ws.Range("L141").Value = 35164.5
$ws.Range("N141").Value = -45524.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 192.73334
$ws.Range("I55").Value = 145.3077
$ws.Range("K55").Value = 145.3077
$ws.Range("M55").Value = 27.69229999999999
$ws.Range("H61").Value = 989.44446
$ws.Range("I61").Value = 873.3333
$ws.Range("J61").Value = 1570
$ws.Range("K61").Value = 873.3333
$ws.Range("L61").Value = 1570
$ws.Range("M61").Value = -671.3333
$ws.Range("N61").Value = -1974
$ws.Range("H113").Value = 989.44446
$ws.Range("I113").Value = 873.3333
$ws.Range("J113").Value = 1570
$ws.Range("K113").Value = 873.3333
$ws.Range("L113").Value = 1570
$ws.Range("M113").Value = 1296.6667
$ws.Range("N113").Value = -5910

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 5858.6665
$ws.Range("J74").Value = 5750.4
$ws.Range("L74").Value = 5750.4
$ws.Range("N74").Value = -7622.4
$ws.Range("H77").Value = 5858.6665
$ws.Range("J77").Value = 5750.4
$ws.Range("L77").Value = 17251.2
$ws.Range("N77").Value = -26611.2
$ws.Range("H107").Value = 643.1667
$ws.Range("I107").Value = 643.1667
$ws.Range("K107").Value = 1929.5001
$ws.Range("M107").Value = -9.500099999999975
$ws.Range("H140").Value = 40571.75
$ws.Range("J140").Value = 40571.75
$ws.Range("L140").Value = 40571.75
$ws.Range("N140").Value = -50931.75
$ws.Range("H141").Value = 45000
$ws.Range("J141").Value = 45000
$ws.Range("L141").Value = 45000
$ws.Range("N141").Value = -55360
